# Add the "Survey 3" child-count figure (row 5, column B) and update the
# active selection, matching the author's manual data entry + the cursor
# position Excel persisted when the file was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 is "Survey 3" (A5); fill in the missing count.
$ws.Range("B5").Value = 23

# Move/restore the selection to where the user left off (C9), matching
# the saved sheetView state.
$ws.Range("C9").Select()
